{"js": "// Update the date heading and every arithmetic answer in the practice\n// table (document was regenerated for a new day: 2025-10-03 Friday).\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph (first paragraph, above the table) ---\nconst firstPara = body.paragraphs.getFirst();\nfirstPara.load(\"text\");\nawait context.sync();\n\nif (firstPara.text.trim() === \"2025-10-02 Thursday\") {\n  firstPara.insertText(\"2025-10-03 Friday\", \"Replace\");\n}\n\n// --- 2. Update every answer cell in the (single) table ---\nconst table = body.tables.getFirst();\ntable.load(\"values\");\nawait context.sync();\n\n// Old expression -> new expression, keyed by exact text, in reading order\n// (row-major, 5 columns per row). A couple of expressions repeat verbatim\n// in the source table, so we replace by position (row/col) rather than by\n// a global text search to avoid touching the wrong occurrence.\nconst newValues = [\n  [\"98-69=29\", \"66-2=64\", \"16+30=46\", \"19+49=68\", \"44-6=38\"],\n  [\"37+43=80\", \"80-29=51\", \"61+18=79\", \"75+10=85\", \"1+90=91\"],\n  [\"37-5=32\", \"78-3=75\", \"92-81=11\", \"18+81=99\", \"18-4=14\"],\n  [\"91+6=97\", \"1+38=39\", \"98-28=70\", \"2+77=79\", \"78-66=12\"],\n  [\"85-73=12\", \"22+48=70\", \"75-55=20\", \"17-8=9\", \"59+17=76\"],\n  [\"2+4=6\", \"51-0=51\", \"36+23=59\", \"37-26=11\", \"37+61=98\"],\n  [\"70+4=74\", \"56+14=70\", \"22+60=82\", \"17+22=39\", \"86-53=33\"],\n  [\"87-12=75\", \"35+51=86\", \"30+7=37\", \"29+15=44\", \"85-27=58\"],\n  [\"16+0=16\", \"88-72=16\", \"65-59=6\", \"68+24=92\", \"26+57=83\"],\n  [\"14-4=10\", \"80+19=99\", \"53-30=23\", \"73-13=60\", \"67+8=75\"],\n  [\"59+36=95\", \"27+52=79\", \"8+43=51\", \"71-56=15\", \"81-23=58\"],\n  [\"50+11=61\", \"20+13=33\", \"74-63=11\", \"46+5=51\", \"61-27=34\"],\n  [\"27+20=47\", \"89-83=6\", \"7+92=99\", \"32+47=79\", \"40-12=28\"],\n  [\"70-37=33\", \"85-11=74\", \"36-3=33\", \"61-17=44\", \"71-34=37\"],\n  [\"54-35=19\", \"48+2=50\", \"39+33=72\", \"56-42=14\", \"52+1=53\"],\n  [\"47+39=86\", \"62+8=70\", \"7+90=97\", \"6+26=32\", \"82-49=33\"],\n  [\"37+10=47\", \"86-35=51\", \"86-41=45\", \"20+54=74\", \"69-43=26\"],\n  [\"24+16=40\", \"64+11=75\", \"27+24=51\", \"76-63=13\", \"7+80=87\"],\n  [\"40-26=14\", \"87-56=31\", \"82-45=37\", \"68-16=52\", \"93-9=84\"],\n  [\"13-5=8\", \"83-67=16\", \"39+52=91\", \"8-4=4\", \"93-39=54\"]\n];\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    table.getCell(r, c).value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and every arithmetic answer in the practice\n# table (document was regenerated for a new day: 2025-10-03 Friday).\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph (first paragraph, above the table) ---\n$p1 = $d.Paragraphs(1)\nif ($p1.Range.Text.TrimEnd(\"`r\") -eq \"2025-10-02 Thursday\") {\n    $p1.Range.Text = \"2025-10-03 Friday\"\n}\n\n# --- 2. Update every answer cell in the (single) table ---\n# Old expression -> new expression, in reading order (row-major, 5 columns\n# per row). A couple of expressions repeat verbatim in the source table, so\n# cells are addressed by (row, col) position rather than by a global\n# text search, which would not disambiguate the duplicate occurrences.\n$newValues = @(\n    @(\"98-69=29\", \"66-2=64\", \"16+30=46\", \"19+49=68\", \"44-6=38\"),\n    @(\"37+43=80\", \"80-29=51\", \"61+18=79\", \"75+10=85\", \"1+90=91\"),\n    @(\"37-5=32\", \"78-3=75\", \"92-81=11\", \"18+81=99\", \"18-4=14\"),\n    @(\"91+6=97\", \"1+38=39\", \"98-28=70\", \"2+77=79\", \"78-66=12\"),\n    @(\"85-73=12\", \"22+48=70\", \"75-55=20\", \"17-8=9\", \"59+17=76\"),\n    @(\"2+4=6\", \"51-0=51\", \"36+23=59\", \"37-26=11\", \"37+61=98\"),\n    @(\"70+4=74\", \"56+14=70\", \"22+60=82\", \"17+22=39\", \"86-53=33\"),\n    @(\"87-12=75\", \"35+51=86\", \"30+7=37\", \"29+15=44\", \"85-27=58\"),\n    @(\"16+0=16\", \"88-72=16\", \"65-59=6\", \"68+24=92\", \"26+57=83\"),\n    @(\"14-4=10\", \"80+19=99\", \"53-30=23\", \"73-13=60\", \"67+8=75\"),\n    @(\"59+36=95\", \"27+52=79\", \"8+43=51\", \"71-56=15\", \"81-23=58\"),\n    @(\"50+11=61\", \"20+13=33\", \"74-63=11\", \"46+5=51\", \"61-27=34\"),\n    @(\"27+20=47\", \"89-83=6\", \"7+92=99\", \"32+47=79\", \"40-12=28\"),\n    @(\"70-37=33\", \"85-11=74\", \"36-3=33\", \"61-17=44\", \"71-34=37\"),\n    @(\"54-35=19\", \"48+2=50\", \"39+33=72\", \"56-42=14\", \"52+1=53\"),\n    @(\"47+39=86\", \"62+8=70\", \"7+90=97\", \"6+26=32\", \"82-49=33\"),\n    @(\"37+10=47\", \"86-35=51\", \"86-41=45\", \"20+54=74\", \"69-43=26\"),\n    @(\"24+16=40\", \"64+11=75\", \"27+24=51\", \"76-63=13\", \"7+80=87\"),\n    @(\"40-26=14\", \"87-56=31\", \"82-45=37\", \"68-16=52\", \"93-9=84\"),\n    @(\"13-5=8\", \"83-67=16\", \"39+52=91\", \"8-4=4\", \"93-39=54\")\n)\n\n$table = $d.Tables(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
